# Update cryptos list (prices & 1h volume % changes); rows 21-22 swap (Uniswap now above BitcoinCash).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price text is purely numeric-looking need NumberFormat forced to
# Text ("@") first, mirroring how these values were authored as text in the source feed
# (Excel would otherwise silently convert e.g. "1.00" -> 1 or "0.270" -> 0.27).
$ws.Range("D2").Value = "79.844.43"
$ws.Range("E2").Value = "  +4.52%  "
$ws.Range("D3").Value = "3.171.59"
$ws.Range("E3").Value = "  +3.15%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "207.49"
$ws.Range("E5").Value = "  +4.66%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "625.39"
$ws.Range("E6").Value = "  +0.91%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.270"
$ws.Range("E7").Value = "  +25.87%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.586"
$ws.Range("E9").Value = "  +6.19%  "
$ws.Range("D10").Value = "3.166.22"
$ws.Range("E10").Value = "  +3.01%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.589"
$ws.Range("E11").Value = "  +32.64%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000254"
$ws.Range("E12").Value = "  +28.81%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.164"
$ws.Range("E13").Value = "  +1.61%  "
$ws.Range("D14").Value = "3.732.71"
$ws.Range("E14").Value = "  +2.47%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.23"
$ws.Range("E15").Value = "  +0.23%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "31.47"
$ws.Range("E16").Value = "  +7.70%  "
$ws.Range("D17").Value = "79.533.43"
$ws.Range("E17").Value = "  +4.27%  "
$ws.Range("D18").Value = "3.146.45"
$ws.Range("E18").Value = "  +1.73%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.22"
$ws.Range("E19").Value = "  +5.42%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.01"
$ws.Range("E20").Value = "  +16.82%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.09"
$ws.Range("E21").Value = "  +0.99%  "
$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "434.16"
$ws.Range("E22").Value = "  +12.86%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.15"
$ws.Range("E23").Value = "  +14.43%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.90"
$ws.Range("E24").Value = "  +6.73%  "
$ws.Range("D25").Value = "3.313.23"
$ws.Range("E25").Value = "  +2.53%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "75.79"
$ws.Range("E26").Value = "  +4.61%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.67"
$ws.Range("E27").Value = "  +2.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.81"
$ws.Range("E28").Value = "  +7.59%  "
$ws.Range("E29").Value = "  -0.23%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0000122"
$ws.Range("E30").Value = "  +12.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  +0.43%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.90"
$ws.Range("E32").Value = "  +7.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "551.12"
$ws.Range("E33").Value = "  +10.04%  "
$ws.Range("E34").Value = "  +2.79%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.151"
$ws.Range("E35").Value = "  +17.16%  "
$ws.Range("E36").Value = "  +2.75%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "22.96"
$ws.Range("E37").Value = "  +10.34%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.122"
$ws.Range("E38").Value = "  +19.40%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.999"
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.404"
$ws.Range("E40").Value = "  +6.71%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "20.75"
$ws.Range("E41").Value = "  +3.43%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "163.29"
$ws.Range("E42").Value = "  -0.06%  "
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.57"
$ws.Range("E44").Value = "  +7.77%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "189.80"
$ws.Range("E45").Value = "  -2.60%  "
$ws.Range("E46").Value = "  +7.83%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.69"
$ws.Range("E47").Value = "  +9.87%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.783"
$ws.Range("E48").Value = "  -2.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.29"
$ws.Range("E49").Value = "  +2.39%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "42.69"
$ws.Range("E50").Value = "  +4.20%  "
$ws.Range("E51").Value = "  +6.76%  "
